$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '64.371.79'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = "'" + '3.136.75'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'" + '608.17'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = "'" + '143.66'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'" + '3.132.99'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = "'" + '0.528'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Value = "'" + '0.150'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').Value = "'" + '5.36'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').Value = "'" + '0.469'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = "'" + '0.0000254'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').Value = "'" + '35.37'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = "'" + '3.654.34'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').Value = "'" + '64.398.02'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = "'" + '3.129.85'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = "'" + '6.86'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = "'" + '476.48'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = "'" + '14.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').Value = "'" + '7.77'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = "'" + '85.42'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.64%  '
$ws.Range('D25').Value = "'" + '13.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -3.35%  '
$ws.Range('D28').Value = "'" + '8.45'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').Value = "'" + '7.27'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.03%  '
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').Value = "'" + '2.05'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.68%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = "'" + '1.00'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'" + '26.75'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.55%  '
$ws.Range('D34').Value = "'" + '2.63'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('D36').Value = "'" + '5.96'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').Value = "'" + '52.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.00%  '
$ws.Range('D38').Value = "'" + '0.0₃0744'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('D39').Value = "'" + '2.99'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('D40').Value = "'" + '445.31'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('D41').Value = "'" + '0.0393'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = "'" + '8.28'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').Value = "'" + '2.889.54'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').Value = "'" + '0.261'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('D46').Value = "'" + '2.23'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').Value = "'" + '2.40'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.63%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = "'" + '0.999'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = "'" + '26.25'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').Value = "'" + '0.113'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = "'" + '120.18'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.68%  '
